$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DiccionarioDatos")

# Clear the descriptive text in B7 (HORA_UTC row) per source data update
$ws.Range("B7").Value = $null

# Reflect the editor's final view state: zoomed in on B7
$ws.Range("B7").Select()
$excel.ActiveWindow.Zoom = 130
